# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 4a0dedb9-... handback row on both the zh-cn and de-de
# report sheets, as part of regenerating the handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-19 10:34:38"
$wsZhCn.Range("H3").Value = "2016-03-19 10:34:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-19 10:34:42"
$wsDeDe.Range("H3").Value = "2016-03-19 10:35:06"
